# Generate Report for Handoff
# - Update "Priority" column (E) to "ht" for the rows whose handoff report
#   was just (re)generated, on both the zh-cn and de-de worksheets.
# - Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for those same rows across all three worksheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 14)

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"          # E: Priority
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-04 22:23:47"  # H: Latest Handoff Datetime
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"          # E: Priority
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-04 22:23:52"  # H: Latest Handoff DateTime
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-04 22:23:52"  # G: Latest HO Xliff Generate Date
}
